# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column K describes the "provincia" field metadata (row2=uri, row3=kind, row4=type).
# It used to be modeled as an sdmx dimension (refArea) with a custom "URI-Provincia"
# datatype; it is now curated as a plain iaest measure of xsd:int.
$ws.Range("K2").Value = "iaest-measure:provincia"
$ws.Range("K3").Value = "medida"
$ws.Range("K4").Value = "xsd:int"
